$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($row, $col, $newText) {
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newText
}

# Row 1 (document row 1)
Set-CellText 1 1 "98÷8="
Set-CellText 1 2 "64÷3="
Set-CellText 1 3 "56÷5="
Set-CellText 1 4 "86÷2="
Set-CellText 1 5 "77÷6="

# Row 5
Set-CellText 5 1 "87÷9="
Set-CellText 5 2 "55÷6="
Set-CellText 5 3 "53÷4="
Set-CellText 5 4 "29÷4="
Set-CellText 5 5 "74÷9="

# Row 9
Set-CellText 9 1 "46÷4="
Set-CellText 9 2 "46÷6="
Set-CellText 9 3 "24÷8="
Set-CellText 9 4 "14÷7="
Set-CellText 9 5 "23÷8="

# Row 13
Set-CellText 13 1 "81÷6="
Set-CellText 13 2 "13÷4="
Set-CellText 13 3 "42÷9="
Set-CellText 13 4 "90÷6="
Set-CellText 13 5 "45÷8="

# Row 17
Set-CellText 17 1 "49÷9="
Set-CellText 17 2 "50÷9="
Set-CellText 17 3 "75÷3="
Set-CellText 17 4 "86÷5="
Set-CellText 17 5 "42÷5="
